$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Columns H/I/J next to the existing red/blue/yellow/green/purple/orange block (rows 2-7) ---
$ws.Range("H2").Value = 1800
$ws.Range("I2").Formula = "=H2/2"

$ws.Range("H3").Value = 1558
$ws.Range("I3:I7").Formula = "=H3/2"
$ws.Range("J3").Formula = "=I3/I2"

$ws.Range("H4").Value = 1272
$ws.Range("J4:J7").Formula = "=I4/I3"

$ws.Range("H5").Value = 1100
$ws.Range("H6").Value = 898
$ws.Range("H7").Value = 634

# --- New row 9 (pixel/round radius block header row) ---
$ws.Range("F9").Value = 378.4
$ws.Range("G9").Formula = "=FLOOR(F9/10.5,2)"
$ws.Range("I9").Value = 900
$ws.Range("J9").Formula = "=I9/25"

# --- Rows 10-14 extend the F/G and I/J blocks ---
$ws.Range("F10").Value = 330
$ws.Range("G10:G13").Formula = "=FLOOR(F10/10.5,2)"
$ws.Range("I10").Value = 775
$ws.Range("J10:J14").Formula = "=I10/25"

$ws.Range("F11").Value = 272.8
$ws.Range("I11").Value = 625

$ws.Range("F12").Value = 238.4
$ws.Range("I12").Value = 550

$ws.Range("F13").Value = 198
$ws.Range("I13").Value = 450

$ws.Range("F14").Value = 145.2
$ws.Range("G14").Formula = "=FLOOR(F14/10.5,2)"
$ws.Range("I14").Value = 325

# --- New rows 16-21 (leg number -> cm block) ---
$ws.Range("B16").Value = 14.14
$ws.Range("F16").Value = 378.4
$ws.Range("G16").Value = 34.57
$ws.Range("H16").Formula = "=F16/G16"

$ws.Range("B17").Value = 12.24
$ws.Range("F17").Value = 330
$ws.Range("G17").Value = 34.57
$ws.Range("H17:H21").Formula = "=F17/G17"
$ws.Range("J17").Formula = "=F17/4.2"

$ws.Range("B18").Value = 10
$ws.Range("F18").Value = 272.8
$ws.Range("G18").Value = 34.57
$ws.Range("J18:J21").Formula = "=F18/4.2"

$ws.Range("B19").Value = 8.66
$ws.Range("F19").Value = 238.4
$ws.Range("G19").Value = 34.57

$ws.Range("B20").Value = 7.07
$ws.Range("F20").Value = 198
$ws.Range("G20").Value = 34.57

$ws.Range("B21").Value = 5
$ws.Range("F21").Value = 145.2
$ws.Range("G21").Value = 34.57

# --- View state: scroll so column B is the leftmost visible column, select G16 ---
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 2
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("G16").Select()

$wb.Save()
